$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.931.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.50%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.570.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.11%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.41%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  -2.19%  "

$ws.Range("E9").Value = "  -2.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.08%  "

$ws.Range("E11").Value = "  -0.34%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.349"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.028.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.813.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000145"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.576.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.61%  "

$ws.Range("E18").Value = "  -3.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "340.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.86%  "

$ws.Range("E20").Value = "  -2.09%  "

$ws.Range("E21").Value = "  -4.08%  "

$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.21%  "

$ws.Range("E26").Value = "  -2.71%  "

$ws.Range("E27").Value = "  -3.84%  "

$ws.Range("E28").Value = "  -2.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "461.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0795"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "176.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.25%  "

$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.398"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.23%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.27%  "

$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("E41").Value = "  -3.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "157.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.93%  "

$ws.Range("E44").Value = "  -3.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.631"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0959"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.22%  "

$ws.Range("E49").Value = "  -1.88%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.12%  "

